$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values that are unambiguous text (contain extra separators / non-numeric chars)
# and therefore remain Text automatically when assigned via .Value
$ws.Range("D2").Value = '35.620.20'
$ws.Range("D3").Value = '1.894.71'
$ws.Range("D15").Value = '2.175.66'
$ws.Range("D17").Value = '1.905.04'
$ws.Range("D18").Value = '35.629.10'
$ws.Range("D20").Value = '0.0₃0831'
$ws.Range("D45").Value = '1.319.15'

# D-column values that look like plain numbers: force Text format first so Excel
# stores them as strings (matching the original inline-string cell type) instead of
# coercing them into numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.59'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.10'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '57.01'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.357'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0753'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0984'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.60'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.797'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.04'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.65'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '246.37'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.02'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.19'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.15'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.95'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.68'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.37'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0608'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.855'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0734'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.24'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.10'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.25'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.35'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.74'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '42.57'

# E-column percentage-change text (always kept as text thanks to the padding spaces)
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -2.17%  '
$ws.Range("E9").Value = '  +9.63%  '
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  +11.23%  '
$ws.Range("E14").Value = '  +9.73%  '
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("E23").Value = '  +4.25%  '
$ws.Range("E24").Value = '  +4.95%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  +4.79%  '
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  +15.76%  '
$ws.Range("E36").Value = '  -16.76%  '
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("E38").Value = '  -3.30%  '
$ws.Range("E39").Value = '  +7.74%  '
$ws.Range("E40").Value = '  +7.23%  '
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("E44").Value = '  +16.74%  '
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("E47").Value = '  +1.73%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("E51").Value = '  -2.54%  '
